$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# Update the "Date" metadata value to reflect the new publication timestamp.
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2022-05-18T17:38:26+00:00"

# Append a new concept row (Level / Code / Display / Definition) for the
# "Targeted Sequencing" experimental strategy.
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A4").Value = "'1"
$concepts.Range("B4").Value = "TARS"
$concepts.Range("C4").Value = "Targeted Sequencing"
$concepts.Range("D4").Value = "Tarteted Sequencing"

# Match the formatting of the preceding data row (style only - the values
# were already set above, and a formats-only paste leaves them untouched).
$concepts.Range("A3:D3").Copy()
$concepts.Range("A4:D4").PasteSpecial($xlPasteFormats)
